# Weekly update: insert a new price observation as row 58, pushing the
# existing rows 58-136 down to 59-137 (the sheet grows from 136 to 137
# data rows, dimension A1:R136 -> A1:R137).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 58; Excel shifts every row at/after
# 58 down by one automatically (formatting such as the date style on
# column D travels with the shift).
$ws.Rows(58).Insert()

# Populate the newly-inserted row 58 with the new observation.
$row = 58
$ws.Cells.Item($row, 1).Value2  = 11                                  # A: Mercado ID
$ws.Cells.Item($row, 2).Value2  = 'Vega Monumental Concepción'        # B: Mercado
$ws.Cells.Item($row, 3).Value2  = 'Bíobío'                            # C: Región
$ws.Cells.Item($row, 4).Value2  = 45174                               # D: Fecha
$ws.Cells.Item($row, 5).Value2  = 8                                   # E: Codreg
$ws.Cells.Item($row, 6).Value2  = 100112037                           # F: Categoría ID
$ws.Cells.Item($row, 7).Value2  = 'Cebollín'                          # G: Categoría
$ws.Cells.Item($row, 8).Value2  = 'Sin especificar'                   # H: Variedad
$ws.Cells.Item($row, 9).Value2  = 'Primera'                           # I: Calidad
$ws.Cells.Item($row, 10).Value2 = 50                                  # J: Volumen
$ws.Cells.Item($row, 11).Value2 = 4500                                # K: Precio mínimo
$ws.Cells.Item($row, 12).Value2 = 4500                                # L: Precio máximo
$ws.Cells.Item($row, 13).Value2 = 4500                                # M: Precio promedio ponderado
$ws.Cells.Item($row, 14).Value2 = '$/paquete 36 unidades'             # N: Unidad de comercialización
$ws.Cells.Item($row, 15).Value2 = 'Región Metropolitana'              # O: Origen
$ws.Cells.Item($row, 16).Value2 = 125                                 # P: Precio $/Kg
$ws.Cells.Item($row, 17).Value2 = 36                                  # Q: Kg o Unidades
$ws.Cells.Item($row, 18).Value2 = 'Hortaliza'                         # R: Clasificación
